$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new 2022-Q3 row at the top
#    of the data block and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 12.23

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 7
$summary.Range("D3").Value = 11.62

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 15
$summary.Range("D4").Value = 14.79

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 23
$summary.Range("D5").Value = 15.73

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 25
$summary.Range("D6").Value = 16

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 18
$summary.Range("D7").Value = 17.57

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 10
$summary.Range("D8").Value = 15.16

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 11
$summary.Range("D9").Value = 6.47

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    sheet that is currently in position 2, "2022-Q2"), and fill it with the
#    fund holding breakdown for that quarter.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "501203"
$q3.Range("C2").Value = "易方达创新未来混合（LOF）"
$q3.Range("D2").Value = "53.16"
$q3.Range("E2").Value = "84.96"
$q3.Range("F2").Value = "7.55"
$q3.Range("G2").Value = "4.0136"
$q3.Range("H2").Value = 3

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "110013"
$q3.Range("C3").Value = "易方达科翔混合"
$q3.Range("D3").Value = "65.50"
$q3.Range("E3").Value = "76.69"
$q3.Range("F3").Value = "4.28"
$q3.Range("G3").Value = "2.8034"
$q3.Range("H3").Value = 3

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "110001"
$q3.Range("C4").Value = "易方达平稳增长混合"
$q3.Range("D4").Value = "34.03"
$q3.Range("E4").Value = "59.68"
$q3.Range("F4").Value = "4.44"
$q3.Range("G4").Value = "1.5109"
$q3.Range("H4").Value = 3

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "009341"
$q3.Range("C5").Value = "易方达均衡成长股票"
$q3.Range("D5").Value = "59.42"
$q3.Range("E5").Value = "87.05"
$q3.Range("F5").Value = "2.42"
$q3.Range("G5").Value = "1.4380"
$q3.Range("H5").Value = 10

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "001018"
$q3.Range("C6").Value = "易方达新经济灵活配置混合"
$q3.Range("D6").Value = "69.74"
$q3.Range("E6").Value = "86.57"
$q3.Range("F6").Value = "1.90"
$q3.Range("G6").Value = "1.3251"
$q3.Range("H6").Value = 10

$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "020010"
$q3.Range("C7").Value = "国泰金牛创新混合"
$q3.Range("D7").Value = "13.26"
$q3.Range("E7").Value = "86.27"
$q3.Range("F7").Value = "3.23"
$q3.Range("G7").Value = "0.4283"
$q3.Range("H7").Value = 8

$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "012173"
$q3.Range("C8").Value = "国泰兴泽优选一年持有期混合A"
$q3.Range("D8").Value = "8.41"
$q3.Range("E8").Value = "88.23"
$q3.Range("F8").Value = "2.50"
$q3.Range("G8").Value = "0.2102"
$q3.Range("H8").Value = 10

$q3.Range("A9").Value = 7
$q3.Range("B9").Value = "012174"
$q3.Range("C9").Value = "国泰兴泽优选一年持有期混合C"
$q3.Range("D9").Value = "6.17"
$q3.Range("E9").Value = "88.23"
$q3.Range("F9").Value = "2.50"
$q3.Range("G9").Value = "0.1542"
$q3.Range("H9").Value = 10

$q3.Range("A10").Value = 8
$q3.Range("B10").Value = "005244"
$q3.Range("C10").Value = "国泰聚优价值灵活配置混合A"
$q3.Range("D10").Value = "4.61"
$q3.Range("E10").Value = "87.30"
$q3.Range("F10").Value = "2.99"
$q3.Range("G10").Value = "0.1378"
$q3.Range("H10").Value = 10

$q3.Range("A11").Value = 9
$q3.Range("B11").Value = "007733"
$q3.Range("C11").Value = "南方智锐混合A"
$q3.Range("D11").Value = "5.65"
$q3.Range("E11").Value = "90.36"
$q3.Range("F11").Value = "2.39"
$q3.Range("G11").Value = "0.1350"
$q3.Range("H11").Value = 6

$q3.Range("A12").Value = 10
$q3.Range("B12").Value = "005245"
$q3.Range("C12").Value = "国泰聚优价值灵活配置混合C"
$q3.Range("D12").Value = "1.80"
$q3.Range("E12").Value = "87.30"
$q3.Range("F12").Value = "2.99"
$q3.Range("G12").Value = "0.0538"
$q3.Range("H12").Value = 10

$q3.Range("A13").Value = 11
$q3.Range("B13").Value = "007734"
$q3.Range("C13").Value = "南方智锐混合C"
$q3.Range("D13").Value = "0.78"
$q3.Range("E13").Value = "90.36"
$q3.Range("F13").Value = "2.39"
$q3.Range("G13").Value = "0.0186"
$q3.Range("H13").Value = 6

# Match the bold/bordered header style used by the sibling quarter sheets.
$headerStyle = $summary.Range("B1").Style
$q3.Range("B1:H1").Style = $headerStyle
$q3.Range("A2:A13").Style = $headerStyle
